# network.xlsx: "Added engine.py fields_map and simulation.py Other changes"
#
# - generators sheet: new "slack" column (V) with header + FALSE flag for
#   every generator row (V2:V7), sheet becomes the active tab with V2:V7
#   selected.
# - crews sheet: no longer the active tab (tabSelected removed as a side
#   effect of activating "generators" instead).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("generators")

# New "slack" header in column V, matching the existing header row style.
$ws.Range("V1").Value = "slack"

# New boolean "slack" flag for each of the 6 generator rows - all FALSE.
$ws.Range("V2:V7").Value = $false

# Make "generators" the active sheet/tab with V2:V7 selected (activeCell
# V2), matching the saved view state in the workbook.
$ws.Activate()
$ws.Range("V2:V7").Select()
